$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add I0 and IF headers, matching style of H1 ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# --- Data rows 2-25: new columns I (I0) and J (IF) ---
$data = @{
    2  = @(1, 4)
    3  = @(7, 8)
    4  = @(1, 1)
    5  = @(3, 7)
    6  = @(1, 4)
    7  = @(1, 5)
    8  = @(1, 6)
    9  = @(2, 7)
    10 = @(1, 5)
    11 = @(1, 4)
    12 = @(7, 8)
    13 = @(5, 7)
    14 = @(7, 9)
    15 = @(5, 8)
    16 = @(5, 7)
    17 = @(8, 9)
    18 = @(8, 9)
    19 = @(7, 8)
    20 = @(7, 9)
    21 = @(7, 9)
    22 = @(6, 7)
    23 = @(4, 6)
    24 = @(5, 6)
    25 = @(3, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
